$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Fill in the "Saker som händer" notes for Vecka 12 (row 9) and Vecka 13 (row 10)
$ws.Range("D9").Value = "Jobbade med rapporten. Skrev klart kod för att både prata med Acc och radion, radion krävde mycket tid."
$ws.Range("D10").Value = "Påskveckan"

# Move the active cell selection from D9 to D10
$ws.Range("D10").Select()
